# Update TPM-recomputed NATMI ligand-receptor metrics (Ccl5-Ccr5) for rows 2-21.
# Values below are taken verbatim from the recomputed dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1780343333333333
$ws.Range("H2").Value = 0.534103
$ws.Range("I2").Value = 0.003649670474736916
$ws.Range("J2").Value = 0.003649670474736915
$ws.Range("M2").Value = 0.243056
$ws.Range("N2").Value = 0.729168
$ws.Range("O2").Value = 0.002199620488481675
$ws.Range("P2").Value = 0.002199620488481675
$ws.Range("Q2").Value = 0.04327231292266667
$ws.Range("R2").Value = 0.389450816304
$ws.Range("S2").Value = [double]"8.027889952437961E-06"
$ws.Range("T2").Value = [double]"8.027889952437959E-06"

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1780343333333333
$ws.Range("H3").Value = 0.534103
$ws.Range("I3").Value = 0.003649670474736916
$ws.Range("J3").Value = 0.003649670474736915
$ws.Range("M3").Value = 70.95253000000001
$ws.Range("N3").Value = 212.85759
$ws.Range("O3").Value = 0.6421097964979703
$ws.Range("P3").Value = 0.6421097964979703
$ws.Range("Q3").Value = 12.63198637686333
$ws.Range("R3").Value = 113.68787739177
$ws.Range("S3").Value = 0.002343489165817972
$ws.Range("T3").Value = 0.002343489165817971

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1780343333333333
$ws.Range("H4").Value = 0.534103
$ws.Range("I4").Value = 0.003649670474736916
$ws.Range("J4").Value = 0.003649670474736915
$ws.Range("M4").Value = 0.04794200000000001
$ws.Range("N4").Value = 0.143826
$ws.Range("O4").Value = 0.0004338679376719292
$ws.Range("P4").Value = 0.0004338679376719292
$ws.Range("Q4").Value = 0.008535322008666668
$ws.Range("R4").Value = 0.076817898078
$ws.Range("S4").Value = [double]"1.583475002056237E-06"
$ws.Range("T4").Value = [double]"1.583475002056237E-06"

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1780343333333333
$ws.Range("H5").Value = 0.534103
$ws.Range("I5").Value = 0.003649670474736916
$ws.Range("J5").Value = 0.003649670474736915
$ws.Range("M5").Value = 39.25553366666666
$ws.Range("N5").Value = 117.766601
$ws.Range("O5").Value = 0.3552567150758761
$ws.Range("P5").Value = 0.3552567150758761
$ws.Range("Q5").Value = 6.988832765989222
$ws.Range("R5").Value = 62.899494893903
$ws.Range("S5").Value = 0.00129656994396445
$ws.Range("T5").Value = 0.00129656994396445

# Row 6
$ws.Range("I6").Value = 0.09908483984804967
$ws.Range("J6").Value = 0.09908483984804965
$ws.Range("M6").Value = 0.243056
$ws.Range("N6").Value = 0.729168
$ws.Range("O6").Value = 0.002199620488481675
$ws.Range("P6").Value = 0.002199620488481675
$ws.Range("Q6").Value = 1.174799266256
$ws.Range("R6").Value = 10.573193396304
$ws.Range("S6").Value = 0.0002179490438276955
$ws.Range("T6").Value = 0.0002179490438276955

# Row 7
$ws.Range("I7").Value = 0.09908483984804967
$ws.Range("J7").Value = 0.09908483984804965
$ws.Range("M7").Value = 70.95253000000001
$ws.Range("N7").Value = 212.85759
$ws.Range("O7").Value = 0.6421097964979703
$ws.Range("P7").Value = 0.6421097964979703
$ws.Range("Q7").Value = 342.9455770810301
$ws.Range("R7").Value = 3086.51019372927
$ws.Range("S7").Value = 0.06362334635086515
$ws.Range("T7").Value = 0.06362334635086514

# Row 8
$ws.Range("I8").Value = 0.09908483984804967
$ws.Range("J8").Value = 0.09908483984804965
$ws.Range("M8").Value = 0.04794200000000001
$ws.Range("N8").Value = 0.143826
$ws.Range("O8").Value = 0.0004338679376719292
$ws.Range("P8").Value = 0.0004338679376719292
$ws.Range("Q8").Value = 0.231725307842
$ws.Range("R8").Value = 2.085527770578
$ws.Range("S8").Value = [double]"4.298973511942671E-05"
$ws.Range("T8").Value = [double]"4.29897351194267E-05"

# Row 9
$ws.Range("I9").Value = 0.09908483984804967
$ws.Range("J9").Value = 0.09908483984804965
$ws.Range("M9").Value = 39.25553366666666
$ws.Range("N9").Value = 117.766601
$ws.Range("O9").Value = 0.3552567150758761
$ws.Range("P9").Value = 0.3552567150758761
$ws.Range("Q9").Value = 189.7396984566837
$ws.Range("R9").Value = 1707.657286110153
$ws.Range("S9").Value = 0.0352005547182374
$ws.Range("T9").Value = 0.0352005547182374

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 35.63223
$ws.Range("H10").Value = 106.89669
$ws.Range("I10").Value = 0.7304540385283456
$ws.Range("J10").Value = 0.7304540385283456
$ws.Range("M10").Value = 0.243056
$ws.Range("N10").Value = 0.729168
$ws.Range("O10").Value = 0.002199620488481675
$ws.Range("P10").Value = 0.002199620488481675
$ws.Range("Q10").Value = 8.660627294880001
$ws.Range("R10").Value = 77.94564565392001
$ws.Range("S10").Value = 0.001606721669041131
$ws.Range("T10").Value = 0.001606721669041131

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 35.63223
$ws.Range("H11").Value = 106.89669
$ws.Range("I11").Value = 0.7304540385283456
$ws.Range("J11").Value = 0.7304540385283456
$ws.Range("M11").Value = 70.95253000000001
$ws.Range("N11").Value = 212.85759
$ws.Range("O11").Value = 0.6421097964979703
$ws.Range("P11").Value = 0.6421097964979703
$ws.Range("Q11").Value = 2528.1968680419
$ws.Range("R11").Value = 22753.7718123771
$ws.Range("S11").Value = 0.4690316940305565
$ws.Range("T11").Value = 0.4690316940305565

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 35.63223
$ws.Range("H12").Value = 106.89669
$ws.Range("I12").Value = 0.7304540385283456
$ws.Range("J12").Value = 0.7304540385283456
$ws.Range("M12").Value = 0.04794200000000001
$ws.Range("N12").Value = 0.143826
$ws.Range("O12").Value = 0.0004338679376719292
$ws.Range("P12").Value = 0.0004338679376719292
$ws.Range("Q12").Value = 1.70828037066
$ws.Range("R12").Value = 15.37452333594
$ws.Range("S12").Value = 0.0003169205872604252
$ws.Range("T12").Value = 0.0003169205872604252

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 35.63223
$ws.Range("H13").Value = 106.89669
$ws.Range("I13").Value = 0.7304540385283456
$ws.Range("J13").Value = 0.7304540385283456
$ws.Range("M13").Value = 39.25553366666666
$ws.Range("N13").Value = 117.766601
$ws.Range("O13").Value = 0.3552567150758761
$ws.Range("P13").Value = 0.3552567150758761
$ws.Range("Q13").Value = 1398.76220438341
$ws.Range("R13").Value = 12588.85983945069
$ws.Range("S13").Value = 0.2594987022414875
$ws.Range("T13").Value = 0.2594987022414875

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.5521946666666667
$ws.Range("H14").Value = 1.656584
$ws.Range("I14").Value = 0.011319887201011
$ws.Range("J14").Value = 0.011319887201011
$ws.Range("M14").Value = 0.243056
$ws.Range("N14").Value = 0.729168
$ws.Range("O14").Value = 0.002199620488481675
$ws.Range("P14").Value = 0.002199620488481675
$ws.Range("Q14").Value = 0.1342142269013334
$ws.Range("R14").Value = 1.207928042112
$ws.Range("S14").Value = [double]"2.489945581464528E-05"
$ws.Range("T14").Value = [double]"2.489945581464528E-05"

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.5521946666666667
$ws.Range("H15").Value = 1.656584
$ws.Range("I15").Value = 0.011319887201011
$ws.Range("J15").Value = 0.011319887201011
$ws.Range("M15").Value = 70.95253000000001
$ws.Range("N15").Value = 212.85759
$ws.Range("O15").Value = 0.6421097964979703
$ws.Range("P15").Value = 0.6421097964979703
$ws.Range("Q15").Value = 39.17960865250667
$ws.Range("R15").Value = 352.61647787256
$ws.Range("S15").Value = 0.007268610467021154
$ws.Range("T15").Value = 0.007268610467021153

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.5521946666666667
$ws.Range("H16").Value = 1.656584
$ws.Range("I16").Value = 0.011319887201011
$ws.Range("J16").Value = 0.011319887201011
$ws.Range("M16").Value = 0.04794200000000001
$ws.Range("N16").Value = 0.143826
$ws.Range("O16").Value = 0.0004338679376719292
$ws.Range("P16").Value = 0.0004338679376719292
$ws.Range("Q16").Value = 0.02647331670933334
$ws.Range("R16").Value = 0.238259850384
$ws.Range("S16").Value = [double]"4.911336114581512E-06"
$ws.Range("T16").Value = [double]"4.911336114581511E-06"

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.5521946666666667
$ws.Range("H17").Value = 1.656584
$ws.Range("I17").Value = 0.011319887201011
$ws.Range("J17").Value = 0.011319887201011
$ws.Range("M17").Value = 39.25553366666666
$ws.Range("N17").Value = 117.766601
$ws.Range("O17").Value = 0.3552567150758761
$ws.Range("P17").Value = 0.3552567150758761
$ws.Range("Q17").Value = 21.67669632788711
$ws.Range("R17").Value = 195.090266950984
$ws.Range("S17").Value = 0.004021465942060622
$ws.Range("T17").Value = 0.004021465942060622

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 7.585023666666667
$ws.Range("H18").Value = 22.755071
$ws.Range("I18").Value = 0.155491563947857
$ws.Range("J18").Value = 0.1554915639478569
$ws.Range("M18").Value = 0.243056
$ws.Range("N18").Value = 0.729168
$ws.Range("O18").Value = 0.002199620488481675
$ws.Range("P18").Value = 0.002199620488481675
$ws.Range("Q18").Value = 1.843585512325334
$ws.Range("R18").Value = 16.592269610928
$ws.Range("S18").Value = 0.0003420224298457647
$ws.Range("T18").Value = 0.0003420224298457646

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 7.585023666666667
$ws.Range("H19").Value = 22.755071
$ws.Range("I19").Value = 0.155491563947857
$ws.Range("J19").Value = 0.1554915639478569
$ws.Range("M19").Value = 70.95253000000001
$ws.Range("N19").Value = 212.85759
$ws.Range("O19").Value = 0.6421097964979703
$ws.Range("P19").Value = 0.6421097964979703
$ws.Range("Q19").Value = 538.1766192598768
$ws.Range("R19").Value = 4843.589573338891
$ws.Range("S19").Value = 0.09984265648370956
$ws.Range("T19").Value = 0.09984265648370955

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 7.585023666666667
$ws.Range("H20").Value = 22.755071
$ws.Range("I20").Value = 0.155491563947857
$ws.Range("J20").Value = 0.1554915639478569
$ws.Range("M20").Value = 0.04794200000000001
$ws.Range("N20").Value = 0.143826
$ws.Range("O20").Value = 0.0004338679376719292
$ws.Range("P20").Value = 0.0004338679376719292
$ws.Range("Q20").Value = 0.3636412046273334
$ws.Range("R20").Value = 3.272770841646
$ws.Range("S20").Value = [double]"6.74628041754396E-05"
$ws.Range("T20").Value = [double]"6.746280417543959E-05"

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 7.585023666666667
$ws.Range("H21").Value = 22.755071
$ws.Range("I21").Value = 0.155491563947857
$ws.Range("J21").Value = 0.1554915639478569
$ws.Range("M21").Value = 39.25553366666666
$ws.Range("N21").Value = 117.766601
$ws.Range("O21").Value = 0.3552567150758761
$ws.Range("P21").Value = 0.3552567150758761
$ws.Range("Q21").Value = 297.7541519092968
$ws.Range("R21").Value = 2679.787367183671
$ws.Range("S21").Value = 0.05523942223012619
$ws.Range("T21").Value = 0.05523942223012619
